$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phen_oncox")

# Update Experimental Factor Ontology (EFO) source_version: v3.62.0 -> v3.63.0
$ws.Range("E4").Value = "v3.63.0"

# Update Disease Ontology source_version: v2024-01-31 -> v2024-02-28
$ws.Range("E3").Value = "v2024-02-28"
